# Apply "repull data, push all data, mean calculation" edits:
# Update column F (dSF) values for several rows to match repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F19").Value = 4
$ws.Range("F21").Value = 16
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 2
$ws.Range("F28").Value = -2
$ws.Range("F29").Value = 2
$ws.Range("F31").Value = 0
$ws.Range("F34").Value = 3
$ws.Range("F36").Value = 2
